$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 74, shifting rows 74:144 down to 75:145.
$ws.Rows.Item(74).Insert()

# Populate the new row 74 with the new weekly record.
$ws.Range("A74").Value2 = 8
$ws.Range("B74").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C74").Value2 = "Coquimbo"
$ws.Range("D74").Value2 = 44705
$ws.Range("E74").Value2 = 4
$ws.Range("F74").Value2 = 100112040
$ws.Range("G74").Value2 = "Cilantro"
$ws.Range("H74").Value2 = "Sin especificar"
$ws.Range("I74").Value2 = "Primera"
$ws.Range("J74").Value2 = 3000
$ws.Range("K74").Value2 = 1500
$ws.Range("L74").Value2 = 2000
$ws.Range("M74").Value2 = 1750
$ws.Range("N74").Value2 = "`$/atado 1 a 1,5 kilos"
$ws.Range("O74").Value2 = "Provincia del Elquí"
$ws.Range("P74").Value2 = 1167
$ws.Range("Q74").Value2 = 1.5
$ws.Range("R74").Value2 = "Hortaliza"
